# Apply cryptos.xlsx price/volume refresh (GitHub Actions data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '28.496.32'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -2.30%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.962.73'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -0.24%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.014'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +0.32%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '322.63'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.72%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '1.011'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.26%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.4785'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -3.88%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.4051'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -4.58%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '54.09'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -0.12%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.08475'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -7.69%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '1.058'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -3.70%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '22.42'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -3.50%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '2.010.54'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +2.59%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '7.592'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -4.27%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '6.167'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -4.47%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '1.014'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +0.41%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '90.53'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.57%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.00001077'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -2.70%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.06626'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -1.14%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '18.58'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -3.82%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '1.012'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +0.39%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '5.871'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -1.25%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '28.564.21'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -2.11%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '11.48'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -4.17%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.301'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.35%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.241.24'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +2.64%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '156.01'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.61%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '20.32'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -1.85%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '5.955'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -4.95%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '2.163'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -4.81%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '124.43'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -2.20%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.9821'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -6.18%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.09643'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -2.27%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.455'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -5.77%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '3.697'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '5.631'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -3.31%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '9.247'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +1.89%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.02334'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -4.23%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.06252'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -1.79%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '1.253'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -2.83%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.6221'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -3.98%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '11.16'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -2.96%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.012'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +0.35%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.1910'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -5.08%  '
$ws.Range("E45").Value = '  +5.80%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '13.12'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -3.17%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.5957'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -4.76%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '2.065'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -5.56%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '3.411'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -2.02%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.06827'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -0.86%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000308'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -9.39%  '
